# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect freshly scraped numbers (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition list) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4468
$ws1.Range("F5").Value = 139
$ws1.Range("F6").Value = 153
$ws1.Range("F10").Value = 621
$ws1.Range("F12").Value = 193
$ws1.Range("F13").Value = 1256
$ws1.Range("F15").Value = 2876
$ws1.Range("F16").Value = 443
$ws1.Range("F17").Value = 576

# --- Sheet "全部类型" (all types, combined list) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4468
$ws4.Range("F5").Value = 139
$ws4.Range("F6").Value = 153
$ws4.Range("F10").Value = 621
$ws4.Range("F13").Value = 193
$ws4.Range("F14").Value = 1256
$ws4.Range("F16").Value = 2876
$ws4.Range("F17").Value = 443
$ws4.Range("F18").Value = 576
